# Update "Percentage of informal employment in total employment" figures
# (total / females / males, columns F/G/H) for the rows whose underlying
# source data was refreshed in this upload.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Tab07")

# Row 50
$ws.Range("F50").Value = 81.599999999999994
$ws.Range("G50").Value = 86.9
$ws.Range("H50").Value = 76

# Row 61
$ws.Range("F61").Value = 87.65
$ws.Range("G61").Value = 90.985714285714295
$ws.Range("H61").Value = 84.8642857142857

# Row 62
$ws.Range("F62").Value = 81.897560975609807
$ws.Range("G62").Value = 83.921951219512195
$ws.Range("H62").Value = 79.914634146341498

# Row 63
$ws.Range("F63").Value = 40.762790697674397
$ws.Range("G63").Value = 40.327906976744202
$ws.Range("H63").Value = 41.0162790697674

# Row 64
$ws.Range("F64").Value = 56.359090909090902
$ws.Range("G64").Value = 54.731818181818198
$ws.Range("H64").Value = 57.572727272727299

# Row 66
$ws.Range("F66").Value = 54.042519685039402
$ws.Range("G66").Value = 54.4015748031496
$ws.Range("H66").Value = 53.574015748031499

# Row 68
$ws.Range("F68").Value = 83.94
$ws.Range("G68").Value = 86.1
$ws.Range("H68").Value = 81.915000000000006

# Row 71
$ws.Range("F71").Value = 87.65
$ws.Range("G71").Value = 90.985714285714295
$ws.Range("H71").Value = 84.8642857142857

# Row 82
$ws.Range("F82").Value = 81.5513513513514
$ws.Range("G82").Value = 83.5324324324325
$ws.Range("H82").Value = 79.583783783783801

# Row 83
$ws.Range("F83").Value = 38.792307692307702
$ws.Range("G83").Value = 38.729487179487201
$ws.Range("H83").Value = 38.9

# Row 84
$ws.Range("F84").Value = 91.0833333333333
$ws.Range("G84").Value = 93.938888888888897
$ws.Range("H84").Value = 88.605555555555597

# Row 90
$ws.Range("F90").Value = 14.0432432432432
$ws.Range("G90").Value = 13.4
$ws.Range("H90").Value = 14.524324324324301

# Row 91
$ws.Range("F91").Value = 89.040740740740802
$ws.Range("G91").Value = 92.248148148148204
$ws.Range("H91").Value = 86.551851851851893

# Row 97
$ws.Range("F97").Value = 86.969230769230805
$ws.Range("G97").Value = 90.838461538461601
$ws.Range("H97").Value = 83.873076923076894
